$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.981.46'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.821.38'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '311.13'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '0.4679'
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('D8').Value = '0.3666'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = '0.07356'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '0.8736'
$ws.Range('D11').Value = '20.30'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '1.827.30'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '5.424'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').Value = '0.07148'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = '6.516'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '91.63'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '0.000008747'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '14.68'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '27.004.13'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = '2.043.50'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').Value = '1.892'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').Value = '150.98'
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').Value = '18.41'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '2.141'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '5.242'
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').Value = '116.70'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').Value = '0.08882'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = '0.7544'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').Value = '1.160'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').Value = '4.503'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D35').Value = '2.944'
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D37').Value = '1.094'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '0.05311'
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').Value = '0.01948'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('D41').Value = '2.379'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '7.188'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').Value = '0.5301'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').Value = '0.1653'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '8.468'
$ws.Range('D46').Value = '0.4898'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('E49').Value = '  -0.99%  '
$ws.Range('D50').Value = '103.09'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').Value = '0.06297'
$ws.Range('E51').Value = '  +0.16%  '
